# Scheduled-runner style data refresh: updates market-price-derived columns
# (H..N: currentAveragePrice*, LevePrice*, LeveProfit*) across the eight
# Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching the
# latest pulled prices. Cells that end up blank are cleared; cells that
# newly have a computed profit are populated.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 41.666668
$ws.Range("I5").Value = 41.666668
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 41.666668
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 73.333332
$ws.Range("N5").Value = $null
# Row 12
$ws.Range("H12").Value = 511.5
$ws.Range("I12").Value = 666
$ws.Range("J12").Value = 460
$ws.Range("K12").Value = 666
$ws.Range("L12").Value = 460
$ws.Range("M12").Value = -496
$ws.Range("N12").Value = -800
# Row 17
$ws.Range("H17").Value = 1025
$ws.Range("J17").Value = 1150
$ws.Range("L17").Value = 3450
$ws.Range("N17").Value = -3786
# Row 92
$ws.Range("H92").Value = 1059.8889
$ws.Range("I92").Value = 1059.8889
$ws.Range("K92").Value = 1059.8889
$ws.Range("M92").Value = 188.1111000000001
# Row 138
$ws.Range("H138").Value = 2247.2222
$ws.Range("I138").Value = 850
$ws.Range("J138").Value = 2421.875
$ws.Range("K138").Value = 2550
$ws.Range("L138").Value = 7265.625
$ws.Range("M138").Value = 2590
$ws.Range("N138").Value = -17545.625

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 5016.625
$ws.Range("I132").Value = 5699.7144
$ws.Range("J132").Value = 235
$ws.Range("K132").Value = 17099.1432
$ws.Range("L132").Value = 705
$ws.Range("M132").Value = -14569.1432
$ws.Range("N132").Value = -5765

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 1413.8572
$ws.Range("J64").Value = 1563.4
$ws.Range("L64").Value = 1563.4
$ws.Range("N64").Value = -2013.4
# Row 67
$ws.Range("H67").Value = 1413.8572
$ws.Range("J67").Value = 1563.4
$ws.Range("L67").Value = 1563.4
$ws.Range("N67").Value = -3123.4
# Row 80
$ws.Range("H80").Value = 1317.5
$ws.Range("I80").Value = 777.5
$ws.Range("J80").Value = 1587.5
$ws.Range("K80").Value = 777.5
$ws.Range("L80").Value = 1587.5
$ws.Range("M80").Value = 220.5
$ws.Range("N80").Value = -3583.5
# Row 83
$ws.Range("H83").Value = 1317.5
$ws.Range("I83").Value = 777.5
$ws.Range("J83").Value = 1587.5
$ws.Range("K83").Value = 3887.5
$ws.Range("L83").Value = 7937.5
$ws.Range("M83").Value = 1104.5
$ws.Range("N83").Value = -17921.5
# Row 94
$ws.Range("H94").Value = 4670
$ws.Range("J94").Value = 4670
$ws.Range("L94").Value = 4670
$ws.Range("N94").Value = -5572
# Row 135
$ws.Range("H135").Value = 67000
$ws.Range("J135").Value = 67000
$ws.Range("L135").Value = 67000
$ws.Range("N135").Value = -77140

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 949.5
$ws.Range("I22").Value = 899
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 899
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -549
$ws.Range("N22").Value = -1700
# Row 31
$ws.Range("H31").Value = 4781.625
$ws.Range("I31").Value = 4278
$ws.Range("K31").Value = 4278
$ws.Range("M31").Value = -3983
# Row 34
$ws.Range("H34").Value = 4781.625
$ws.Range("I34").Value = 4278
$ws.Range("K34").Value = 4278
$ws.Range("M34").Value = -4076
# Row 122
$ws.Range("H122").Value = 866.3333
$ws.Range("I122").Value = 866.3333
$ws.Range("K122").Value = 2598.9999
$ws.Range("M122").Value = -148.9998999999998
# Row 134
$ws.Range("H134").Value = 2500
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = $null
$ws.Range("N134").Value = -12570

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 87.333336
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null
# Row 92
$ws.Range("H92").Value = 95
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
# Row 122
$ws.Range("H122").Value = 827.5
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 770
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 6930
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -11830

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7000
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730
# Row 73
$ws.Range("H73").Value = 7000
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 58
$ws.Range("H58").Value = 49875
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 49875
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 49875
$ws.Range("M58").Value = $null
$ws.Range("N58").Value = -50395
# Row 61
$ws.Range("H61").Value = 1255
$ws.Range("I61").Value = 885
$ws.Range("J61").Value = 1995
$ws.Range("K61").Value = 885
$ws.Range("L61").Value = 1995
$ws.Range("M61").Value = -683
$ws.Range("N61").Value = -2399
# Row 113
$ws.Range("H113").Value = 1255
$ws.Range("I113").Value = 885
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 885
$ws.Range("L113").Value = 1995
$ws.Range("M113").Value = 1285
$ws.Range("N113").Value = -6335
# Row 122
$ws.Range("H122").Value = 15192.75
$ws.Range("I122").Value = 15192.75
$ws.Range("K122").Value = 45578.25
$ws.Range("M122").Value = -43128.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = $null
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
# Row 45
$ws.Range("H45").Value = 29999
$ws.Range("J45").Value = 29999
$ws.Range("L45").Value = 29999
$ws.Range("N45").Value = -30981
# Row 51
$ws.Range("H51").Value = 48794.2
$ws.Range("I51").Value = 48794.2
$ws.Range("K51").Value = 48794.2
$ws.Range("M51").Value = -48284.2
# Row 96
$ws.Range("H96").Value = 3330
$ws.Range("I96").Value = 3250
$ws.Range("K96").Value = 3250
$ws.Range("M96").Value = -1877
# Row 100
$ws.Range("H100").Value = 500
$ws.Range("I100").Value = 366.66666
$ws.Range("K100").Value = 733.33332
$ws.Range("M100").Value = -192.33332
# Row 105
$ws.Range("H105").Value = 57950
$ws.Range("J105").Value = 57950
$ws.Range("L105").Value = 57950
$ws.Range("N105").Value = -64938
# Row 107
$ws.Range("H107").Value = 1613.875
$ws.Range("I107").Value = 1018
$ws.Range("K107").Value = 3054
$ws.Range("M107").Value = -1134
